# Scheduled-runner market refresh: update per-item price/profit
# columns (H:N) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve
# tables to the latest Universalis snapshot.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H74").Value = 4773.4443
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H77").Value = 4773.4443
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H92").Value = 1262.3334
$ws.Range("I92").Value = 1031.4736
$ws.Range("J92").Value = 2139.6
$ws.Range("K92").Value = 1031.4736
$ws.Range("L92").Value = 2139.6
$ws.Range("M92").Value = 216.5264
$ws.Range("N92").Value = -4635.6
$ws.Range("H100").Value = 2097.2856
$ws.Range("I100").Value = 1450
$ws.Range("K100").Value = 1450
$ws.Range("M100").Value = -909
$ws.Range("H113").Value = 14288220
$ws.Range("I113").Value = 18184190
$ws.Range("J113").Value = 2995
$ws.Range("K113").Value = 18184190
$ws.Range("L113").Value = 2995
$ws.Range("M113").Value = -18180936
$ws.Range("N113").Value = -9503
$ws.Range("H129").Value = 872.12195
$ws.Range("J129").Value = 883.7105
$ws.Range("L129").Value = 2651.1315
$ws.Range("N129").Value = -12651.1315
$ws.Range("H138").Value = 2448.4316
$ws.Range("J138").Value = 2687.3975
$ws.Range("L138").Value = 8062.1925
$ws.Range("N138").Value = -18342.1925

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 16000000
$ws.Range("I13").Value = 16000000
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 16000000
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -15999856
$ws.Range("N13").ClearContents()
$ws.Range("H32").Value = 9810.484
$ws.Range("I32").Value = 7191.2783
$ws.Range("J32").Value = 20156.35
$ws.Range("K32").Value = 7191.2783
$ws.Range("L32").Value = 20156.35
$ws.Range("M32").Value = -6904.2783
$ws.Range("N32").Value = -20730.35
$ws.Range("H45").Value = 1263.5
$ws.Range("I45").Value = 1011.125
$ws.Range("K45").Value = 1011.125
$ws.Range("M45").Value = -634.125
$ws.Range("H74").Value = 1955.6333
$ws.Range("I74").Value = 1031.3334
$ws.Range("J74").Value = 4112.3335
$ws.Range("K74").Value = 1031.3334
$ws.Range("L74").Value = 4112.3335
$ws.Range("M74").Value = -157.3334
$ws.Range("N74").Value = -5860.3335
$ws.Range("H77").Value = 1955.6333
$ws.Range("I77").Value = 1031.3334
$ws.Range("J77").Value = 4112.3335
$ws.Range("K77").Value = 5156.666999999999
$ws.Range("L77").Value = 20561.6675
$ws.Range("M77").Value = -788.6669999999995
$ws.Range("N77").Value = -29297.6675
$ws.Range("H110").Value = 931.1111
$ws.Range("I110").Value = 925.7143
$ws.Range("K110").Value = 925.7143
$ws.Range("M110").Value = 1119.2857
$ws.Range("H132").Value = 3557.1177
$ws.Range("I132").Value = 3189.2222
$ws.Range("J132").Value = 3971
$ws.Range("K132").Value = 9567.6666
$ws.Range("L132").Value = 11913
$ws.Range("M132").Value = -7037.6666
$ws.Range("N132").Value = -16973

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 90
$ws.Range("I7").Value = 90
$ws.Range("K7").Value = 90
$ws.Range("M7").Value = 23
$ws.Range("H35").Value = 29800
$ws.Range("J35").Value = 29800
$ws.Range("L35").Value = 29800
$ws.Range("N35").Value = -30420
$ws.Range("H107").Value = 1082.6666
$ws.Range("I107").Value = 1017.4545
$ws.Range("K107").Value = 1017.4545
$ws.Range("M107").Value = 902.5454999999999
$ws.Range("H134").Value = 6212.8945
$ws.Range("I134").Value = 863.4375
$ws.Range("K134").Value = 2590.3125
$ws.Range("M134").Value = -55.3125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 668.8333
$ws.Range("I105").Value = 656.9091
$ws.Range("K105").Value = 656.9091
$ws.Range("M105").Value = 1090.0909
$ws.Range("H107").Value = 943.9524
$ws.Range("I107").Value = 602.375
$ws.Range("J107").Value = 2037
$ws.Range("K107").Value = 602.375
$ws.Range("L107").Value = 2037
$ws.Range("M107").Value = 1317.625
$ws.Range("N107").Value = -5877
$ws.Range("H111").Value = 44500
$ws.Range("J111").Value = 44500
$ws.Range("L111").Value = 44500
$ws.Range("N111").Value = -52680
$ws.Range("H134").Value = 13890341
$ws.Range("I134").Value = 1479.5
$ws.Range("J134").Value = 31251418
$ws.Range("K134").Value = 4438.5
$ws.Range("L134").Value = 93754254
$ws.Range("M134").Value = -1903.5
$ws.Range("N134").Value = -93759324

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 202400.8
$ws.Range("H107").Value = 9716.182000000001
$ws.Range("J107").Value = 14928.286
$ws.Range("L107").Value = 44784.858
$ws.Range("N107").Value = -48624.858
$ws.Range("H129").Value = 23149334
$ws.Range("I129").Value = 66667092
$ws.Range("J129").Value = 6411734.5
$ws.Range("K129").Value = 200001276
$ws.Range("L129").Value = 19235203.5
$ws.Range("M129").Value = -199996276
$ws.Range("N129").Value = -19245203.5
$ws.Range("H130").Value = 2570.5625
$ws.Range("J130").Value = 2570.5625
$ws.Range("L130").Value = 7711.6875
$ws.Range("N130").Value = -17751.6875
$ws.Range("H140").Value = 32659.363
$ws.Range("I140").Value = 51978.95
$ws.Range("J140").Value = 2936.923
$ws.Range("K140").Value = 155936.85
$ws.Range("L140").Value = 8810.769
$ws.Range("M140").Value = -150756.85
$ws.Range("N140").Value = -19170.769

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1164.9286
$ws.Range("I97").Value = 1070.8
$ws.Range("K97").Value = 1070.8
$ws.Range("M97").Value = -574.8
$ws.Range("H107").Value = 3846734
$ws.Range("I107").Value = 4808317.5
$ws.Range("J107").Value = 400
$ws.Range("K107").Value = 4808317.5
$ws.Range("L107").Value = 400
$ws.Range("M107").Value = -4806397.5
$ws.Range("N107").Value = -4240
$ws.Range("H126").Value = 1715.9333
$ws.Range("I126").Value = 1518.3846
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 4555.1538
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -2085.1538
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 12271
$ws.Range("I132").Value = 30652.75
$ws.Range("J132").Value = 4101.3335
$ws.Range("K132").Value = 91958.25
$ws.Range("L132").Value = 12304.0005
$ws.Range("M132").Value = -89428.25
$ws.Range("N132").Value = -17364.0005

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5500
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5500
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5500
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -5876

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H124").Value = 20666.666
$ws.Range("J124").Value = 20666.666
$ws.Range("L124").Value = 20666.666
$ws.Range("N124").Value = -30486.666
$ws.Range("H126").Value = 37038132
$ws.Range("I126").Value = 52911000
$ws.Range("J126").Value = 1442.3334
$ws.Range("K126").Value = 158733000
$ws.Range("L126").Value = 4327.0002
$ws.Range("M126").Value = -158730530
$ws.Range("N126").Value = -9267.0002
$ws.Range("H132").Value = 2736.4312
$ws.Range("I132").Value = 2777.8125
$ws.Range("K132").Value = 8333.4375
$ws.Range("M132").Value = -5803.4375
